$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 used to hold the "Total Duration:" / "28 Hours" summary text in
# C22/D22. Replace it with a new timesheet entry: a date in A22 and a
# clock-in time in B22 (stored as literal text, matching the other date /
# time columns in this sheet), and clear the old summary text that used
# to live in C22/D22 so they become blank cells again.

# Force text interpretation (so "2026-02-07" / "12:30:38" aren't
# auto-converted into date/time serial numbers) by pre-setting the
# number format to Text before assigning the values.
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "2026-02-07"

$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "12:30:38"

# Restore the original cell formatting (style index 2, same as the rest
# of the sheet's date/time text cells) by copying the format from the
# neighboring C22/D22 cells, which already carry that style.
$ws.Range("C22").Copy()
$ws.Range("A22").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("D22").Copy()
$ws.Range("B22").PasteSpecial(-4122)  # xlPasteFormats

# Clear the old "Total Duration:" / "28 Hours" text, leaving blank cells
# with the same style as before.
$ws.Range("C22").ClearContents()
$ws.Range("D22").ClearContents()

$excel.CutCopyMode = $false
